# Participant 3, session 1, screen condition: headset disconnected around
# the end of part 1 and had to be restarted, so the raw log had to be
# stitched back together by hand. This adds the "Year since course"
# column (new col C) and fills in the part-1/part-2 rows for participants
# 2-4 (par id column B / new col C) plus the missing dates for rows 4-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "date" column (C), shifting the
# old C:F ("date","condition","date","condition") right to D:G.
$ws.Columns("C:C").Insert()

# par id column (B): fill in the condition codes for participants 2-4.
$ws.Range("B3").Value = "SY"
$ws.Range("B4").Value = "NZ"

# New column header + "years since course" values for participants 1-3.
$ws.Range("C1").Value = "Year since course"
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 5

$ws.Range("B5").Value = "DT"

# Participants 3 and 4 ran on day two (session restarted after the
# headset disconnect), so their date column is filled in separately.
$ws.Range("D4").Value = 20250108
$ws.Range("D5").Value = 20250108

# Column C width follows Excel's own best-fit sizing for the new header.
$ws.Columns("C:C").ColumnWidth = 14.3

# Leave the selection where the editor last left off.
$ws.Range("C5").Select() | Out-Null
